$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they stay text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "29.157.03"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").Value = "1.848.55"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "0.7027"
$ws.Range("E5").Value = "  -4.87%  "
$ws.Range("D6").Value = "238.76"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.3059"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").Value = "0.07496"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("D11").Value = "0.08125"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").Value = "1.857.21"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("D13").Value = "0.7261"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "5.228"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "88.70"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("D16").Value = "29.262.05"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "5.765"
$ws.Range("E17").Value = "  -6.34%  "
$ws.Range("D18").Value = "238.34"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("D20").Value = "0.000007631"
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "2.117.27"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "7.597"
$ws.Range("E24").Value = "  -4.11%  "
$ws.Range("D25").Value = "9.001"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").Value = "161.14"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "0.1454"
$ws.Range("E27").Value = "  -7.42%  "
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").Value = "1.993"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").Value = "1.396"
$ws.Range("E30").Value = "  -5.54%  "
$ws.Range("D31").Value = "4.550"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "1.492"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "3.973"
$ws.Range("E33").Value = "  -5.41%  "
$ws.Range("D34").Value = "0.05167"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "1.186"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "1.039"
$ws.Range("E36").Value = "  +3.67%  "
$ws.Range("D37").Value = "0.7023"
$ws.Range("E37").Value = "  -8.83%  "
$ws.Range("D38").Value = "2.650"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").Value = "0.01865"
$ws.Range("E39").Value = "  -4.92%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Value = "0.9341"
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "1.076.82"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").Value = "0.4284"
$ws.Range("D45").Value = "70.16"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "102.65"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("D48").Value = "2.011.13"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "1.743"
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.047"
$ws.Range("E50").Value = "  -7.23%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.153"
$ws.Range("E51").Value = "  -4.67%  "
